$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four oldest years (2000, 2002, 2005, 2007). The remaining
# rows (2010, 2012, 2015, 2017) shift up to become rows 2-5.
$ws.Rows("2:5").Delete()

# Copy the formatting of the year-label column down onto the new row
# so the 2020 row matches the styling (bold, centered, bordered) used
# by the other year cells in column A.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# Append the new 2020 row of data.
$ws.Range("A6").Value = "2020年"
$ws.Range("B6").Value = 4554576.12250768
$ws.Range("C6").Value = 841651090.814226
$ws.Range("D6").Value = 330800455.825737
$ws.Range("F6").Value = 3840083149.58612
$ws.Range("J6").Value = 276123836.16857
$ws.Range("K6").Value = 126886993.017075
$ws.Range("L6").Value = 116095264.790195
$ws.Range("M6").Value = 131207526.795881
$ws.Range("O6").Value = 1222279.57958158
$ws.Range("P6").Value = 4880487.07462051
$ws.Range("R6").Value = 5486655.16723377
$ws.Range("S6").Value = 608515165.754948
